$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 6
$ws.Range("B4").Value = 55
$ws.Range("C6").Value = 11
$ws.Range("B7").Value = 56
$ws.Range("C7").Value = 19
$ws.Range("B8").Value = 42
$ws.Range("C8").Value = 23
$ws.Range("B9").Value = 48
$ws.Range("C10").Value = 17
$ws.Range("B11").Value = 39
$ws.Range("C12").Value = 21
$ws.Range("B13").Value = 21
$ws.Range("C13").Value = 12
